# Updated cryptos list on Tue Dec 26 09:44:38 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to remain a text/string cell (preserving exact formatting,
    # e.g. trailing zeros like "1.01" or "231.40") instead of letting Excel
    # auto-convert parseable numeric-looking strings into real numbers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Column D (Price) updates ---
Set-TextValue $ws.Range("D2")  "42.839.84"
Set-TextValue $ws.Range("D3")  "2.236.58"
Set-TextValue $ws.Range("D4")  "1.01"
Set-TextValue $ws.Range("D5")  "114.52"
Set-TextValue $ws.Range("D6")  "275.28"
Set-TextValue $ws.Range("D10") "46.42"
Set-TextValue $ws.Range("D14") "15.25"
Set-TextValue $ws.Range("D16") "2.575.06"
Set-TextValue $ws.Range("D17") "2.244.03"
Set-TextValue $ws.Range("D18") "42.809.23"
Set-TextValue $ws.Range("D20") "6.75"
Set-TextValue $ws.Range("D21") "72.11"
Set-TextValue $ws.Range("D24") "231.40"
Set-TextValue $ws.Range("D29") "40.11"
Set-TextValue $ws.Range("D32") "173.10"
Set-TextValue $ws.Range("D35") "5.56"
Set-TextValue $ws.Range("D36") "4.39"
Set-TextValue $ws.Range("D39") "0.0370"
Set-TextValue $ws.Range("D42") "71.25"
Set-TextValue $ws.Range("D46") "1.33"
Set-TextValue $ws.Range("D49") "8.44"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value  = "  -0.96%  "
$ws.Range("E3").Value  = "  -1.49%  "
$ws.Range("E4").Value  = "  +0.26%  "
$ws.Range("E5").Value  = "  +3.35%  "
$ws.Range("E6").Value  = "  +4.63%  "
$ws.Range("E7").Value  = "  -2.97%  "
$ws.Range("E8").Value  = "  +0.01%  "
$ws.Range("E9").Value  = "  +0.91%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  -3.68%  "
$ws.Range("E23").Value = "  +5.09%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  +5.89%  "
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  +12.73%  "
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("E42").Value = "  -5.52%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -0.21%  "

# --- Row 43: Algorand -> Celestia ---
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D43") "13.22"
$ws.Range("E43").Value = "  -7.18%  "

# --- Row 44: Celestia -> Algorand ---
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D44") "0.233"
$ws.Range("E44").Value = "  -0.98%  "

# --- Row 51: Aave -> TheSandbox ---
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D51") "0.642"
$ws.Range("E51").Value = "  +8.32%  "
